$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.458714962005615
$ws.Range("B1").Value = 1.788643956184387
$ws.Range("C1").Value = 1.701697468757629
$ws.Range("D1").Value = 1.567126393318176
$ws.Range("E1").Value = 1.101616501808167
